# Insert a new weekly price record at row 518 of the single data sheet.
# This pushes the previous rows 518:543 down to 519:544 (dimension grows
# from A1:R543 to A1:R544) and fills the newly-opened row 518 with the
# latest "Acelga" (Macroferia Regional de Talca) observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row above the current row 518, shifting 518:543 -> 519:544.
$ws.Rows.Item(518).Insert()

# Populate the new row 518 with the new weekly observation.
$ws.Cells.Item(518, 1).Value = 5
$ws.Cells.Item(518, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(518, 3).Value = "Maule"
$ws.Cells.Item(518, 4).Value = 45267
$ws.Cells.Item(518, 5).Value = 7
$ws.Cells.Item(518, 6).Value = 100112009
$ws.Cells.Item(518, 7).Value = "Acelga"
$ws.Cells.Item(518, 8).Value = "Sin especificar"
$ws.Cells.Item(518, 9).Value = "Primera"
$ws.Cells.Item(518, 10).Value = 500
$ws.Cells.Item(518, 11).Value = 2000
$ws.Cells.Item(518, 12).Value = 2000
$ws.Cells.Item(518, 13).Value = 2000
$ws.Cells.Item(518, 14).Value = "$/docena de atados (4 kilos)"
$ws.Cells.Item(518, 15).Value = "Región del Maule"
$ws.Cells.Item(518, 16).Value = 500
$ws.Cells.Item(518, 17).Value = 4
$ws.Cells.Item(518, 18).Value = "Hortaliza"
